$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "Våler (Østfold)"
$ws.Range("B10").Value = "Indre Østfold"
$ws.Range("B17").Value = "Nordre Follo"
$ws.Range("B41").Value = "Øvre Eiker"
$ws.Range("B53").Value = "Nore og Uvdal"
$ws.Range("B66").Value = "Våler (Innlandet)"
$ws.Range("B89").Value = "Østre Toten"
$ws.Range("B90").Value = "Vestre Toten"
$ws.Range("B92").Value = "Søndre Land"
$ws.Range("B93").Value = "Nordre Land"
$ws.Range("B97").Value = "Vestre Slidre"
$ws.Range("B98").Value = "Øystre Slidre"
$ws.Range("B138").Value = "Evje og Hornnes"
$ws.Range("B219").Value = "Herøy (Møre og Romsdal)"
$ws.Range("B241").Value = "Trondheim - Tråante"
$ws.Range("B243").Value = "Namsos - Nåavmesjenjaelmie"
$ws.Range("B248").Value = "Røros - Rosse"
$ws.Range("B250").Value = "Midtre Gauldal"
$ws.Range("B259").Value = "Levanger - Levangke"
$ws.Range("B261").Value = "Snåase - Snåsa"
$ws.Range("B263").Value = "Raarvihke - Røyrvik"
$ws.Range("B271").Value = "Indre Fosen"
$ws.Range("B286").Value = "Herøy (Nordland)"
$ws.Range("B291").Value = "Aarborte - Hattfjelldal"
$ws.Range("B295").Value = "Rana - Raane"
$ws.Range("B303").Value = "Fauske - Fuossko"
$ws.Range("B304").Value = "Sørfold - Fuolldá"
$ws.Range("B307").Value = "Evenes - Evenássi"
$ws.Range("B314").Value = "Bø (Nordland)"
$ws.Range("B316").Value = "Sortland - Suortá"
$ws.Range("B319").Value = "Hábmer - Hamarøy"
$ws.Range("B321").Value = "Harstad - Hárstták"
$ws.Range("B323").Value = "Dielddanuorri - Tjeldsund"
$ws.Range("B325").Value = "Gratangen - Rivtták"
$ws.Range("B326").Value = "Loabák - Lavangen"
$ws.Range("B336").Value = "Storfjord - Omasvuotna - Omasvuono"
$ws.Range("B337").Value = "Gáivuotna - Kåfjord - Kaivuono"
$ws.Range("B339").Value = "Nordreisa - Ráisa - Raisi"
$ws.Range("B342").Value = "Hammerfest - Hámmerfeasta"
$ws.Range("B345").Value = "Kárásjohka - Karasjok"
$ws.Range("B346").Value = "Guovdageaidnu - Kautokeino"
$ws.Range("B351").Value = "Porsanger - Porsángu - Porsanki"
$ws.Range("B354").Value = "Deatnu - Tana"
$ws.Range("B358").Value = "Unjárga - Nesseby"
